# The edit splits two existing sentences so that "etc" and "deep sea" sit in
# their own runs bracketed by <w:proofErr> spell/grammar markers, and appends
# three new outline paragraphs (plus one blank paragraph) after the second
# paragraph, moving the trailing _GoBack bookmark onto the final paragraph.
#
# We rebuild the whole body in one shot via Range.InsertXML so paragraph
# indices/bookmarks stay consistent (re-querying Paragraphs(n) between
# separate InsertXML calls is unreliable once the body has been mutated).

$d = $word.ActiveDocument

$bodyXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">Categorize samples into MAB, SNE, GOM </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>etc</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> regions. Tally the most abundant fishes in each region.</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t>Categorize fishes into deep sea and shallow.</w:t></w:r>
<w:r><w:t xml:space="preserve"> Calculate the percentage of </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>deep sea</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> species at surface depth vs at deep depth.</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t xml:space="preserve">Compare </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>eDNA</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> with BTS</w:t></w:r>
</w:p>
<w:p/>
<w:p>
<w:r><w:t>Deep vs demersal vs pelagic vs mesopelagic (took all surface vs &gt;200m samples and run indicator species analysis)</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t>Vertical structure: thermocline vs chlorophyll max, which matters to fish</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

[void]$d.Content.InsertXML($bodyXml)
